$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 45.85166310918853

# Row 3
$ws.Range("B3").Value = 0.127881588408715
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 35.318310910425

# Row 4
$ws.Range("B4").Value = 0.127881588408715
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 1.094976487407548
